# chore: update Sheets via scheduled runner
# Refresh market-price-derived columns (currentAveragePrice*, LevePrice*, LeveProfit*)
# across the per-job profit tables with the latest pulled values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 210762.81
$ws.Range("J17").Value = 234309.89
$ws.Range("L17").Value = 702929.67
$ws.Range("N17").Value = -703265.67

$ws.Range("H94").Value = 6481.769
$ws.Range("I94").Value = 1033.125
$ws.Range("J94").Value = 15199.6
$ws.Range("K94").Value = 1033.125
$ws.Range("L94").Value = 15199.6
$ws.Range("M94").Value = -582.125
$ws.Range("N94").Value = -16101.6

$ws.Range("H100").Value = 5919.4165
$ws.Range("I100").Value = 6003
$ws.Range("K100").Value = 6003
$ws.Range("M100").Value = -5462

$ws.Range("H112").Value = 2119.3157
$ws.Range("J112").Value = 1886.2941
$ws.Range("L112").Value = 5658.8823
$ws.Range("N112").Value = -7874.8823

$ws.Range("H138").Value = 8122.2104
$ws.Range("I138").Value = 7561.143
$ws.Range("J138").Value = 8248.903
$ws.Range("K138").Value = 22683.429
$ws.Range("L138").Value = 24746.709
$ws.Range("M138").Value = -17543.429
$ws.Range("N138").Value = -35026.709

$ws.Range("H141").Value = 2836
$ws.Range("I141").Value = 2836
$ws.Range("K141").Value = 8508
$ws.Range("M141").Value = -3328

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1266.15
$ws.Range("I2").Value = 958.7857
$ws.Range("K2").Value = 958.7857
$ws.Range("M2").Value = -845.7857

$ws.Range("H32").Value = 198701.86
$ws.Range("I32").Value = 235019.23
$ws.Range("J32").Value = 3496
$ws.Range("K32").Value = 235019.23
$ws.Range("L32").Value = 3496
$ws.Range("M32").Value = -234732.23
$ws.Range("N32").Value = -4070

$ws.Range("H45").Value = 399175.84
$ws.Range("I45").Value = 619327.75
$ws.Range("J45").Value = 2902.4
$ws.Range("K45").Value = 619327.75
$ws.Range("L45").Value = 2902.4
$ws.Range("M45").Value = -618950.75
$ws.Range("N45").Value = -3656.4

$ws.Range("H110").Value = 35721510
$ws.Range("I110").Value = 43479828
$ws.Range("K110").Value = 43479828
$ws.Range("M110").Value = -43477783

$ws.Range("H116").Value = 1266.15
$ws.Range("I116").Value = 958.7857
$ws.Range("K116").Value = 958.7857
$ws.Range("M116").Value = 1335.2143

$ws.Range("H132").Value = 25642218
$ws.Range("I132").Value = 26316962
$ws.Range("K132").Value = 78950886
$ws.Range("M132").Value = -78948356

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1266.15
$ws.Range("I3").Value = 958.7857
$ws.Range("K3").Value = 958.7857
$ws.Range("M3").Value = -844.7857

$ws.Range("H99").Value = 702.6842
$ws.Range("J99").Value = 636.3333
$ws.Range("L99").Value = 636.3333
$ws.Range("N99").Value = -3632.3333

$ws.Range("H105").Value = 1444.3043
$ws.Range("I105").Value = 1372.6666
$ws.Range("J105").Value = 1702.2
$ws.Range("K105").Value = 1372.6666
$ws.Range("L105").Value = 1702.2
$ws.Range("M105").Value = 374.3334
$ws.Range("N105").Value = -5196.2

$ws.Range("H107").Value = 83417290
$ws.Range("I107").Value = 34583.332
$ws.Range("K107").Value = 34583.332
$ws.Range("M107").Value = -32663.332

$ws.Range("H134").Value = 3915.75
$ws.Range("I134").Value = 3915.75
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 11747.25
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -9212.25
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1255.5294
$ws.Range("I16").Value = 1259.5555
$ws.Range("K16").Value = 1259.5555
$ws.Range("M16").Value = -972.5554999999999

$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()

$ws.Range("H31").Value = 2532.5806
$ws.Range("I31").Value = 955.7857
$ws.Range("K31").Value = 955.7857
$ws.Range("M31").Value = -660.7857

$ws.Range("H34").Value = 2532.5806
$ws.Range("I34").Value = 955.7857
$ws.Range("K34").Value = 955.7857
$ws.Range("M34").Value = -753.7857

$ws.Range("H58").Value = 296309.25
$ws.Range("I58").Value = 1075.6111
$ws.Range("K58").Value = 1075.6111
$ws.Range("M58").Value = -872.6111000000001

$ws.Range("H105").Value = 12895.8
$ws.Range("I105").Value = 12895.8
$ws.Range("K105").Value = 12895.8
$ws.Range("M105").Value = -11148.8

$ws.Range("H107").Value = 4026.6667
$ws.Range("I107").Value = 4026.6667
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 4026.6667
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -2106.6667
$ws.Range("N107").ClearContents()

$ws.Range("H113").Value = 1255.5294
$ws.Range("I113").Value = 1259.5555
$ws.Range("K113").Value = 1259.5555
$ws.Range("M113").Value = 910.4445000000001

$ws.Range("H136").Value = 296309.25
$ws.Range("I136").Value = 1075.6111
$ws.Range("K136").Value = 3226.8333
$ws.Range("M136").Value = -676.8333000000002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H31").Value = 1124.75
$ws.Range("I31").Value = 1124.75
$ws.Range("K31").Value = 3374.25
$ws.Range("M31").Value = -3086.25

$ws.Range("H68").Value = 2021.1364
$ws.Range("I68").Value = 1264.75
$ws.Range("K68").Value = 3794.25
$ws.Range("M68").Value = -2983.25

$ws.Range("H71").Value = 2021.1364
$ws.Range("I71").Value = 1264.75
$ws.Range("K71").Value = 11382.75
$ws.Range("M71").Value = -7326.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 1004241.6
$ws.Range("I3").Value = 200
$ws.Range("J3").Value = 1255252
$ws.Range("K3").Value = 200
$ws.Range("L3").Value = 1255252
$ws.Range("M3").Value = -84
$ws.Range("N3").Value = -1255484

$ws.Range("H70").Value = 7476.6206
$ws.Range("I70").Value = 7335.95
$ws.Range("K70").Value = 7335.95
$ws.Range("M70").Value = -7065.95

$ws.Range("H73").Value = 7476.6206
$ws.Range("I73").Value = 7335.95
$ws.Range("K73").Value = 7335.95
$ws.Range("M73").Value = -6399.95

$ws.Range("H113").Value = 2331
$ws.Range("I113").Value = 1538.1111
$ws.Range("J113").Value = 3223
$ws.Range("K113").Value = 1538.1111
$ws.Range("L113").Value = 3223
$ws.Range("M113").Value = 631.8888999999999
$ws.Range("N113").Value = -7563

$ws.Range("H126").Value = 6483.16
$ws.Range("I126").Value = 7581.1177
$ws.Range("K126").Value = 22743.3531
$ws.Range("M126").Value = -20273.3531

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1919.9722
$ws.Range("J82").Value = 3021.2727
$ws.Range("L82").Value = 3021.2727
$ws.Range("N82").Value = -3743.2727

$ws.Range("H85").Value = 1919.9722
$ws.Range("J85").Value = 3021.2727
$ws.Range("L85").Value = 3021.2727
$ws.Range("N85").Value = -5517.2727

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 2502
$ws.Range("I2").Value = 2502
$ws.Range("K2").Value = 2502
$ws.Range("M2").Value = -2390

$ws.Range("H106").Value = 24585
$ws.Range("L106").Value = 24585
$ws.Range("N106").Value = -27109

$ws.Range("H107").Value = 888.2
$ws.Range("I107").Value = 599
$ws.Range("J107").Value = 960.5
$ws.Range("K107").Value = 1797
$ws.Range("L107").Value = 2881.5
$ws.Range("M107").Value = 123
$ws.Range("N107").Value = -6721.5

$ws.Range("H113").Value = 2778.3125
$ws.Range("J113").Value = 5664.3335
$ws.Range("L113").Value = 16993.0005
$ws.Range("N113").Value = -21333.0005
